{"js": "// The document edit updates the date line and every \"A\u00d7B=C\" multiplication\n// answer in the table to a new value. Every old value in this document is\n// unique, so a direct text search + replace (one hit each) reproduces the\n// diff exactly, regardless of which table cell/paragraph it lives in.\nconst pairs = [\n  [\"2024-09-19 Thursday\", \"2024-09-20 Friday\"],\n  [\"703\u00d72=1406\", \"989\u00d78=7912\"],\n  [\"908\u00d74=3632\", \"997\u00d77=6979\"],\n  [\"873\u00d77=6111\", \"528\u00d74=2112\"],\n  [\"888\u00d74=3552\", \"681\u00d75=3405\"],\n  [\"338\u00d73=1014\", \"240\u00d76=1440\"],\n  [\"430\u00d76=2580\", \"543\u00d75=2715\"],\n  [\"120\u00d72=240\", \"947\u00d72=1894\"],\n  [\"795\u00d75=3975\", \"918\u00d75=4590\"],\n  [\"539\u00d73=1617\", \"408\u00d79=3672\"],\n  [\"101\u00d78=808\", \"357\u00d78=2856\"],\n  [\"923\u00d73=2769\", \"929\u00d77=6503\"],\n  [\"609\u00d73=1827\", \"361\u00d72=722\"],\n  [\"265\u00d76=1590\", \"356\u00d72=712\"],\n  [\"938\u00d73=2814\", \"612\u00d79=5508\"],\n  [\"560\u00d73=1680\", \"319\u00d72=638\"],\n  [\"129\u00d75=645\", \"285\u00d75=1425\"],\n  [\"509\u00d74=2036\", \"635\u00d75=3175\"],\n  [\"549\u00d74=2196\", \"973\u00d77=6811\"],\n  [\"481\u00d78=3848\", \"413\u00d72=826\"],\n  [\"356\u00d76=2136\", \"833\u00d77=5831\"],\n  [\"290\u00d73=870\", \"583\u00d73=1749\"],\n  [\"456\u00d73=1368\", \"232\u00d76=1392\"],\n  [\"394\u00d77=2758\", \"121\u00d78=968\"],\n  [\"196\u00d72=392\", \"830\u00d74=3320\"],\n  [\"721\u00d73=2163\", \"861\u00d75=4305\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  // matchCase keeps the search exact (e.g. \"240\" shouldn't match inside\n  // \"1240\"-style numbers, and none of these values are substrings of\n  // one another anyway).\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document edit updates the date line and every \"A\u00d7B=C\" multiplication\n# answer in the table to a new value. Every old value in this document is\n# unique, so Find/Replace with MatchCase = true and Replace = wdReplaceOne\n# (one hit each) reproduces the diff exactly, regardless of which table\n# cell/paragraph it lives in.\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceOne = 1\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\nfunction Replace-OneMatch($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #         ReplaceWith, Replace)\n    $ok = $find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceOne)\n    if (-not $ok) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-OneMatch '2024-09-19 Thursday' '2024-09-20 Friday'\nReplace-OneMatch '703\u00d72=1406' '989\u00d78=7912'\nReplace-OneMatch '908\u00d74=3632' '997\u00d77=6979'\nReplace-OneMatch '873\u00d77=6111' '528\u00d74=2112'\nReplace-OneMatch '888\u00d74=3552' '681\u00d75=3405'\nReplace-OneMatch '338\u00d73=1014' '240\u00d76=1440'\nReplace-OneMatch '430\u00d76=2580' '543\u00d75=2715'\nReplace-OneMatch '120\u00d72=240' '947\u00d72=1894'\nReplace-OneMatch '795\u00d75=3975' '918\u00d75=4590'\nReplace-OneMatch '539\u00d73=1617' '408\u00d79=3672'\nReplace-OneMatch '101\u00d78=808' '357\u00d78=2856'\nReplace-OneMatch '923\u00d73=2769' '929\u00d77=6503'\nReplace-OneMatch '609\u00d73=1827' '361\u00d72=722'\nReplace-OneMatch '265\u00d76=1590' '356\u00d72=712'\nReplace-OneMatch '938\u00d73=2814' '612\u00d79=5508'\nReplace-OneMatch '560\u00d73=1680' '319\u00d72=638'\nReplace-OneMatch '129\u00d75=645' '285\u00d75=1425'\nReplace-OneMatch '509\u00d74=2036' '635\u00d75=3175'\nReplace-OneMatch '549\u00d74=2196' '973\u00d77=6811'\nReplace-OneMatch '481\u00d78=3848' '413\u00d72=826'\nReplace-OneMatch '356\u00d76=2136' '833\u00d77=5831'\nReplace-OneMatch '290\u00d73=870' '583\u00d73=1749'\nReplace-OneMatch '456\u00d73=1368' '232\u00d76=1392'\nReplace-OneMatch '394\u00d77=2758' '121\u00d78=968'\nReplace-OneMatch '196\u00d72=392' '830\u00d74=3320'\nReplace-OneMatch '721\u00d73=2163' '861\u00d75=4305'\n"}
